# Change utilization choice "in_store_for_allocation" to "not_installed"
# and update the active sheet / selection state.

$wb = $excel.ActiveWorkbook

# --- Data change -----------------------------------------------------
# "choices" sheet, row 15 is the current_use choice that used to be
# in_store_for_allocation / In Store For Allocation / Almacenado
# Epsperando asignación. Re-purpose it as not_installed.
$choices = $wb.Worksheets.Item("choices")
$choices.Range("B15").Value = "not_installed"
$choices.Range("C15").Value = "Not Installed"
$choices.Range("D15").Value = "No Instalado"

# --- View / selection changes -----------------------------------------
# The "choices" sheet used to be the active tab with the whole row 12
# selected; now it is just a background sheet with a single cell (C12)
# selected, and "survey" becomes the active tab instead.
$choices.Activate()
$choices.Range("C12").Select()

$survey = $wb.Worksheets.Item("survey")
$survey.Activate()
